# Refitting NCDEs to individual patients (for manuscript figure)
#
# Adds a new "Label" column (H) to Sheet1 marking each patient row as
# Control (0) or MDD (1), and refreshes a handful of re-fit metric values
# (D/E columns) that shifted slightly after the NCDEs were refit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Label" header in H1, matching the style of the other headers ---
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Label"

# --- Updated metric values (D/E/F) from the re-fit ---
$ws.Range("D3").Value = 0.3174436629119137
$ws.Range("E3").Value = 0.3174436629119137

$ws.Range("D9").Value = 0.366992896129867
$ws.Range("E9").Value = 0.633007103870133

$ws.Range("D10").Value = 0.1203821889845419
$ws.Range("E10").Value = 0.8796178110154581

$ws.Range("F11").Value = 1.040018796920776

# --- New "Label" column values: 0 = Control, 1 = MDD ---
$labels = @{
    2 = 0; 3 = 0; 4 = 0; 5 = 0; 6 = 0;
    7 = 1; 8 = 1; 9 = 1; 10 = 1; 11 = 1;
    12 = 0; 13 = 0; 14 = 0; 15 = 0; 16 = 0;
    17 = 1; 18 = 1; 19 = 1; 20 = 1; 21 = 1;
}

foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 8).Value = $labels[$row]
}
